$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the two new rows with shared string values
$ws.Cells.Item(12, 2).Value = "Minor"
$ws.Cells.Item(13, 2).Value = "Major"

# Apply the Minor/Major theme font scheme to the new cells' fonts
$ws.Cells.Item(12, 2).Font.ThemeFont = 1
$ws.Cells.Item(13, 2).Font.ThemeFont = 2
